$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @{
    2 = "model_10_5_0"
    3 = "model_10_5_22"
    4 = "model_10_5_21"
    5 = "model_10_5_20"
    6 = "model_10_5_19"
    7 = "model_10_5_18"
    8 = "model_10_5_17"
    9 = "model_10_5_16"
    10 = "model_10_5_15"
    11 = "model_10_5_14"
    12 = "model_10_5_13"
    13 = "model_10_5_23"
    14 = "model_10_5_12"
    15 = "model_10_5_10"
    16 = "model_10_5_9"
    17 = "model_10_5_8"
    18 = "model_10_5_7"
    19 = "model_10_5_6"
    20 = "model_10_5_5"
    21 = "model_10_5_4"
    22 = "model_10_5_3"
    23 = "model_10_5_2"
    24 = "model_10_5_1"
    25 = "model_10_5_11"
    26 = "model_10_5_24"
}

$rowValues = @(0.6731329884640765, 0.4382949588972336, 0.8554970179463681, 0.5115830679343389, 0.3617455065250397, 0.7840085029602051, 0.01952212303876877, 0.4242500960826874)

foreach ($r in 2..26) {
    $ws.Cells.Item($r, 1).Value = $names[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $rowValues[$c]
    }
}
